# "fixed hip profile issue"
#
# Adds a new PAN-OS "set" command line to the "set commands" sheet:
#   set rulebase security rules ""Outbound Block Rule"" source-hip any
#
# The canonical OOXML diff shows this new line being inserted twice into
# the sheet's single data column (once in the "Outbound Block Rule" block,
# right after the "...service any" line, and a second time lower down in
# the "Inbound Block Rule" block, right after its own "...service any"
# line) which pushes every row below each insertion point down by one.
# Reproduce that here with two ordinary row inserts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set commands")

$newLine = 'set rulebase security rules ""Outbound Block Rule"" source-hip any'

# First insertion: directly above the existing
#   set rulebase security rules ""Outbound Block Rule"" action deny
# row (currently row 629), i.e. right after "...service any" for the
# Outbound Block Rule.
$ws.Range("A629").EntireRow.Insert()
$ws.Range("A629").Value = $newLine

# Second insertion: directly above the existing
#   set rulebase security rules ""Inbound Block Rule"" action deny
# row. Because of the insert above, that row is now at 642 (was 641).
$ws.Range("A642").EntireRow.Insert()
$ws.Range("A642").Value = $newLine
